$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 415; existing rows 415..447 shift down to 416..448
$ws.Rows.Item(415).Insert()

# Populate the newly inserted row 415 with the new weekly data point
$ws.Cells.Item(415, 1).Value = 9
$ws.Cells.Item(415, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(415, 3).Value = "Metropolitana"
$ws.Cells.Item(415, 4).Value = 45166
$ws.Cells.Item(415, 5).Value = 13
$ws.Cells.Item(415, 6).Value = 100112043
$ws.Cells.Item(415, 7).Value = "Pepino ensalada"
$ws.Cells.Item(415, 8).Value = "Sin especificar"
$ws.Cells.Item(415, 9).Value = "Primera"
$ws.Cells.Item(415, 10).Value = 70
$ws.Cells.Item(415, 11).Value = 11000
$ws.Cells.Item(415, 12).Value = 12000
$ws.Cells.Item(415, 13).Value = 11500
$ws.Cells.Item(415, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(415, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(415, 16).Value = 192
$ws.Cells.Item(415, 17).Value = 60
$ws.Cells.Item(415, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date-time style used by the other date cells in column D
$ws.Cells.Item(415, 4).NumberFormat = $ws.Cells.Item(416, 4).NumberFormat
